$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the sheet's inlineStr/shared-string cells)
# instead of auto-converting them to numeric values.
$textCells = @("D5", "D8", "D15", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D32", "D34", "D36", "D40", "D42", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.452.73"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.565.08"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "208.39"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D8").Value = "22.08"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "1.786.32"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "1.564.28"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "63.56"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "27.431.89"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "212.93"
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "9.58"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").Value = "152.84"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D27").Value = "6.67"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "14.97"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "1.373.78"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "2.97"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").Value = "0.821"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D42").Value = "0.978"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D44").Value = "63.97"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").Value = "1.698.80"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "85.42"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "0.0957"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "0.0496"
$ws.Range("E51").Value = "  -0.55%  "

# Restore default styling on the cells we temporarily formatted as text,
# so no stray number format remains applied.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
